# edit.ps1 - reproduces the OOXML diff:
#   1. (best effort) empty p15:sldGuideLst extension on the presentation
#   2. datetimeFigureOut placeholder text "2019. 6. 5." -> "2019. 9. 28."
#      on the slide master, every slide layout, and the notes master
#   3. slide 5 "tensor size" paragraph: "[3 x 5 + 80]" -> "[3 x (5 + 80)]"
#      (typed over the "x 5 + 80]" selection, producing two runs)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Touch the presentation-level Guides collection. Real PowerPoint
#    leaves behind an empty <p:extLst><p:ext .../><p15:sldGuideLst/>
#    ...</p:ext></p:extLst> after the guides UI has been touched even
#    when no guide is actually added. Best effort only - harmless if
#    unsupported.
# ---------------------------------------------------------------------
try {
    $app = $ppt
    $app.DisplayGuides = $true
} catch {
}
try {
    $null = $p.Guides.Add(1, 100)
} catch {
}

# ---------------------------------------------------------------------
# 2) Refresh the auto date field text everywhere it is cached:
#    slide master, every custom layout, and the notes master.
# ---------------------------------------------------------------------
$oldDate = "2019. 6. 5."
$newDate = "2019. 9. 28."

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

try {
    if ($p.HasNotesMaster) {
        $notesMaster = $p.NotesMaster
        Update-DateShapes $notesMaster.Shapes
    } else {
        $notesMaster = $p.NotesMaster
        Update-DateShapes $notesMaster.Shapes
    }
} catch {
}

# ---------------------------------------------------------------------
# 3) Slide 5: "... tensor size : N x N x [3 x 5 + 80]"
#            -> "... tensor size : N x N x [3 x (5 + 80)]"
#    Done by retyping the "x 5 + 80]" selection as "x (5 + 80)]",
#    which is exactly what the author's edit did (splits the run in two).
# ---------------------------------------------------------------------
$needle = "x 5 + 80]"
$replacement = "x (5 + 80)]"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                $full = $tr.Text
                $idx0 = $full.IndexOf($needle)
                if ($idx0 -ge 0) {
                    $start = $idx0 + 1
                    $len = $needle.Length
                    $sub = $tr.Characters($start, $len)
                    $sub.Text = $replacement
                }
            }
        }
    }
}
